$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.452571110594402
$ws.Range("D2").Value = 0.6552889247863019

$ws.Range("C3").Value = -0.1316788246857113
$ws.Range("D3").Value = 0.8964351649087434

$ws.Range("C4").Value = -1.305017093907319
$ws.Range("D4").Value = 0.2053726424877964

$ws.Range("C5").Value = -0.2868839979017052
$ws.Range("D5").Value = 0.7768857158919804

$ws.Range("C6").Value = -0.7041369195877211
$ws.Range("D6").Value = 0.4887344825019164

$ws.Range("C7").Value = -1.530976341307013
$ws.Range("D7").Value = 0.1400290789488416

$ws.Range("C8").Value = -0.7221292309079389
$ws.Range("D8").Value = 0.4778212511897366

$ws.Range("C9").Value = -1.341311078092691
$ws.Range("D9").Value = 0.1935032232562928

$ws.Range("C10").Value = -0.2014115297430123
$ws.Range("D10").Value = 0.8422274323289975

$ws.Range("C11").Value = 1.23918500581923
$ws.Range("D11").Value = 0.2283347970114582
